# Weekly data refresh: insert 4 new price rows (week of 2022-01-21, serial 44582)
# at the top of the existing Ciruela records (previously starting at row 84),
# pushing all prior rows down by 4 (84-129 -> 88-133).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows before the current row 84 (shifts rows 84:129 -> 88:133).
$ws.Range("A84:A87").EntireRow.Insert() | Out-Null

# Values shared by every record in this block.
$mercadoId = 4
$mercado   = "Feria Lagunitas de Puerto Montt"
$region    = "Los Lagos"
$codreg    = 10
$tipo      = "Fruta"
$productoId = 100103
$producto   = "Frutos de hueso (carozo)"
$categoriaId = 100103002
$categoria   = "Ciruela"

$fecha = 44582

$nuevasFilas = @(
    @{ Row=84; Variedad="Black Amber"; Calidad="Primera"; Volumen=400; PMin=16000; PMax=16500; PProm=16250; Unidad="`$/caja 15 kilos granel"; Origen="Región de O'Higgins"; PKg=1083; KgUnidad=15 },
    @{ Row=85; Variedad="Black Amber"; Calidad="Segunda"; Volumen=200; PMin=14000; PMax=14000; PProm=14000; Unidad="`$/caja 15 kilos granel"; Origen="Región de O'Higgins"; PKg=933;  KgUnidad=15 },
    @{ Row=86; Variedad="Lemon";       Calidad="Primera"; Volumen=400; PMin=16000; PMax=16500; PProm=16250; Unidad="`$/caja 15 kilos granel"; Origen="Región de O'Higgins"; PKg=1083; KgUnidad=15 },
    @{ Row=87; Variedad="Lemon";       Calidad="Segunda"; Volumen=200; PMin=14000; PMax=14000; PProm=14000; Unidad="`$/caja 15 kilos granel"; Origen="Región de O'Higgins"; PKg=933;  KgUnidad=15 }
)

foreach ($fila in $nuevasFilas) {
    $r = $fila.Row

    $ws.Cells.Item($r, 1).Value  = $mercadoId
    $ws.Cells.Item($r, 2).Value  = $mercado
    $ws.Cells.Item($r, 3).Value  = $region
    $ws.Cells.Item($r, 4).Value  = $fecha
    $ws.Cells.Item($r, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r, 5).Value  = $codreg
    $ws.Cells.Item($r, 6).Value  = $tipo
    $ws.Cells.Item($r, 7).Value  = $productoId
    $ws.Cells.Item($r, 8).Value  = $producto
    $ws.Cells.Item($r, 9).Value  = $categoriaId
    $ws.Cells.Item($r, 10).Value = $categoria
    $ws.Cells.Item($r, 11).Value = $fila.Variedad
    $ws.Cells.Item($r, 12).Value = $fila.Calidad
    $ws.Cells.Item($r, 13).Value = $fila.Volumen
    $ws.Cells.Item($r, 14).Value = $fila.PMin
    $ws.Cells.Item($r, 15).Value = $fila.PMax
    $ws.Cells.Item($r, 16).Value = $fila.PProm
    $ws.Cells.Item($r, 17).Value = $fila.Unidad
    $ws.Cells.Item($r, 18).Value = $fila.Origen
    $ws.Cells.Item($r, 19).Value = $fila.PKg
    $ws.Cells.Item($r, 20).Value = $fila.KgUnidad
}
